# Applies the timetable regeneration fix described in the commit message.
# Updates Section_A, Section_B and Course_Summary sheets with the new
# course codes / names / tutorial slots / instructor assignments.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Section_A
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "Free"
$wsA.Range("C2").Value = "CS312"
$wsA.Range("D2").Value = "CS307 (Tutorial)"
$wsA.Range("F2").Value = "CS307"

$wsA.Range("B3").Value = "CS312 (Tutorial)"
$wsA.Range("C3").Value = "Free"
$wsA.Range("E3").Value = "CS312"
$wsA.Range("F3").Value = "CS308"

$wsA.Range("C5").Value = "CS307"
$wsA.Range("D5").Value = "CS312"
$wsA.Range("E5").Value = "Free"
$wsA.Range("F5").Value = "CS308 (Tutorial)"

$wsA.Range("B6").Value = "Free"
$wsA.Range("D6").Value = "CS308"
$wsA.Range("E6").Value = "CS465 (Elective)"
$wsA.Range("F6").Value = "Free"

$wsA.Range("B7").Value = "Free"
$wsA.Range("C7").Value = "CS308"
$wsA.Range("D7").Value = "CS307"
$wsA.Range("E7").Value = "Free"

# ---------------------------------------------------------------------
# Section_B
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("C2").Value = "Free"
$wsB.Range("E2").Value = "CS307"
$wsB.Range("F2").Value = "CS312"

$wsB.Range("B3").Value = "CS307 (Tutorial)"
$wsB.Range("D3").Value = "CS307"
$wsB.Range("E3").Value = "CS312"
$wsB.Range("F3").Value = "Free"

$wsB.Range("B5").Value = "CS312"
$wsB.Range("C5").Value = "CS308 (Tutorial)"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "Free"
$wsB.Range("F5").Value = "Free"

$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "CS308"
$wsB.Range("E6").Value = "CS465 (Elective)"

$wsB.Range("B7").Value = "CS308"
$wsB.Range("C7").Value = "CS307"
$wsB.Range("D7").Value = "CS312 (Tutorial)"
$wsB.Range("E7").Value = "CS308"

# ---------------------------------------------------------------------
# Course_Summary
# ---------------------------------------------------------------------
$wsC = $wb.Worksheets.Item("Course_Summary")

$wsC.Range("A2").Value = "CS312"
$wsC.Range("B2").Value = "Data Analysis for CS"
$wsC.Range("F2").Value = "Dr. Rohit Singh"

$wsC.Range("A3").Value = "CS307"
$wsC.Range("B3").Value = "Advanced Networks"
$wsC.Range("F3").Value = "Dr. Pooja Nair"

$wsC.Range("A4").Value = "CS308"
$wsC.Range("B4").Value = "Machine Learning"
$wsC.Range("F4").Value = "Dr. Karthik Rao"

$wsC.Range("A5").Value = "CS465"
$wsC.Range("B5").Value = "Distributed Systems"
$wsC.Range("F5").Value = "Dr. Divya Gupta"
